$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AYKO")

$row = 64

$values = @(
    "6178",
    "6/18/2025",
    "LA PAMPA 5368",
    "12",
    "807658629",
    "AYKO",
    "Pendiente",
    "Poste inclinado",
    "1",
    "Aplomo",
    "Sin equipos",
    "Poste"
)

for ($col = 1; $col -le 12; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$col - 1]
    $cell.ClearFormats()
}

$ws.Cells.Item($row, 13).Value = -58.482752
$ws.Cells.Item($row, 14).Value = -34.581371
